$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet holds a DFA transition table in A1:E7 (header in row 1,
# six data rows in rows 2-7). The same six-row block is appended three more
# times so the table spans rows 8-25.
$block = @(
    @("q0", "1", "q1", $true,  $true),
    @("q0", "0", "q0", $true,  $true),
    @("q1", "1", "q2", $false, $true),
    @("q1", "0", "q0", $false, $true),
    @("q2", "1", "q2", $false, $false),
    @("q2", "0", "q2", $false, $false)
)

$row = 8
for ($rep = 0; $rep -lt 3; $rep++) {
    foreach ($r in $block) {
        $ws.Cells.Item($row, 1).Value = $r[0]
        # Column B holds the digit as text ("1"/"0"), not a number, to match
        # the source data - the leading apostrophe forces text entry; the
        # style is then reset so no quote-prefix formatting carries over.
        $ws.Cells.Item($row, 2).Value = "'" + $r[1]
        $ws.Cells.Item($row, 2).Style = "Normal"
        $ws.Cells.Item($row, 3).Value = $r[2]
        $ws.Cells.Item($row, 4).Value = $r[3]
        $ws.Cells.Item($row, 5).Value = $r[4]
        $row++
    }
}
